$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the two existing 2-row merges in column A; the dataset groups
#     are growing from 2 rows to 3 rows each, so the merge layout changes.
$ws.Range("A2:A3").UnMerge()
$ws.Range("A5:A6").UnMerge()

# --- Prepare formatting for the six brand-new rows (8-13) by copying the
#     existing A:B cell formatting (bold / border / centered style) down.
$ws.Range("A2:B2").Copy()
$ws.Range("A8:B13").PasteSpecial(-4122) # xlPasteFormats

# --- Row 2 (enzymes, iter 1)
$ws.Range("A2").Value = "enzymes"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 35.5
$ws.Range("D2").Value = 0.0298328677803526
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 3

# --- Row 3 (enzymes, iter 2)
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 35.66666666666667
$ws.Range("D3").Value = 0.0469515116310907
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 2

# --- Row 4 (enzymes, iter 3) -- new row
$ws.Range("A4").ClearContents()
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 35.83333333333334
$ws.Range("D4").Value = 0.027588242262078
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1

# --- Row 5 (imdb, iter 1)
$ws.Range("A5").Value = "imdb"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 72.89999999999999
$ws.Range("D5").Value = 0.0244049175372505
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 2

# --- Row 6 (imdb, iter 2)
$ws.Range("A6").ClearContents()
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 70.7
$ws.Range("D6").Value = 0.0194010309004444
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 2

# --- Row 7 (imdb, iter 3) -- new row
$ws.Range("A7").ClearContents()
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 72
$ws.Range("D7").Value = 0.0375233260785874
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1

# --- Row 8 (mutag, iter 1) -- new row
$ws.Range("A8").Value = "mutag"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 82.99999999999999
$ws.Range("D8").Value = 0.08390470785361211
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1

# --- Row 9 (mutag, iter 2) -- new row
$ws.Range("A9").ClearContents()
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 81.5
$ws.Range("D9").Value = 0.0490917508345343
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 2

# --- Row 10 (mutag, iter 3) -- new row
$ws.Range("A10").ClearContents()
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 76.5
$ws.Range("D10").Value = 0.0548634668973808
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 3

# --- Row 11 (proteins, iter 1) -- new row
$ws.Range("A11").Value = "proteins"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 70.98214285714286
$ws.Range("D11").Value = 0.0247757802241278
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 1

# --- Row 12 (proteins, iter 2) -- new row
$ws.Range("A12").ClearContents()
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 72.58928571428571
$ws.Range("D12").Value = 0.0236969062821908
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 3

# --- Row 13 (proteins, iter 3) -- new row
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 72.14285714285714
$ws.Range("D13").Value = 0.0225594452819339
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 2

# --- Re-merge column A into 3-row dataset groups.
$ws.Range("A2:A4").Merge()
$ws.Range("A5:A7").Merge()
$ws.Range("A8:A10").Merge()
$ws.Range("A11:A13").Merge()

# Merging auto-generates "split border" formatting (top/middle/bottom
# segments) for the newly merged ranges; restore the original uniform
# style (same bordered/bold/centered look used throughout column A) so
# every A cell keeps the same style index it had before merging.
$ws.Range("B2").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A5:A7").PasteSpecial(-4122)
$ws.Range("A8:A10").PasteSpecial(-4122)
$ws.Range("A11:A13").PasteSpecial(-4122)
